# Update "想去人数" (interest count) figures on the 展览 and 全部类型 sheets
# to reflect the latest scrape (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1880
$ws1.Range("F8").Value  = 23
$ws1.Range("F11").Value = 105
$ws1.Range("F13").Value = 771
$ws1.Range("F14").Value = 214
$ws1.Range("F18").Value = 339
$ws1.Range("F19").Value = 195
$ws1.Range("F20").Value = 688
$ws1.Range("F21").Value = 62
$ws1.Range("F25").Value = 896
$ws1.Range("F26").Value = 339
$ws1.Range("F29").Value = 296

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 1880
$ws4.Range("F10").Value = 23
$ws4.Range("F13").Value = 105
$ws4.Range("F15").Value = 771
$ws4.Range("F16").Value = 214
$ws4.Range("F23").Value = 339
$ws4.Range("F27").Value = 195
$ws4.Range("F28").Value = 688
$ws4.Range("F29").Value = 62
$ws4.Range("F33").Value = 896
$ws4.Range("F34").Value = 339
$ws4.Range("F39").Value = 296
